$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold, border, centered) from A4 to A5, matching the year-label column style
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new 2021 data row
$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = 12585
$ws.Range("C5").Value = 3676
$ws.Range("D5").Value = 882
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 9216
$ws.Range("G5").Value = 172443
$ws.Range("H5").Value = 10429
$ws.Range("I5").Value = 802
$ws.Range("J5").Value = 5042
$ws.Range("K5").Value = 2617
$ws.Range("L5").Value = 2900
$ws.Range("M5").Value = 178553
$ws.Range("N5").Value = 794
$ws.Range("O5").Value = 95
$ws.Range("P5").Value = 3896
$ws.Range("Q5").Value = 3038
$ws.Range("R5").Value = 275
$ws.Range("S5").Value = 2690
$ws.Range("T5").Value = 9385
$ws.Range("U5").Value = 664
$ws.Range("V5").Value = 7706
$ws.Range("W5").Value = 77
$ws.Range("X5").Value = 889
$ws.Range("Y5").Value = 818
$ws.Range("Z5").Value = 3749
$ws.Range("AA5").Value = 2267
$ws.Range("AB5").Value = 15371
$ws.Range("AC5").Value = 2275
$ws.Range("AD5").Value = 822
$ws.Range("AE5").Value = 49
$ws.Range("AF5").Value = 6267
$ws.Range("AG5").Value = 3486
$ws.Range("AH5").Value = 13457
$ws.Range("AI5").Value = 14626
$ws.Range("AJ5").Value = 2614
$ws.Range("AK5").Value = 2765
$ws.Range("AL5").Value = 2361
$ws.Range("AM5").Value = 209
$ws.Range("AN5").Value = 12111
$ws.Range("AO5").Value = 2552
$ws.Range("AP5").Value = 13897
$ws.Range("AQ5").Value = 779
$ws.Range("AR5").Value = 4509
$ws.Range("AS5").Value = 1747
$ws.Range("AT5").Value = 272

Write-Host "Row 5 (2021 data) added successfully"
